{"js": "// Apply Garrett's Findings 1.25.23 meeting-notes update:\n//  - Paragraph 1 (\"A\") becomes the PC.skel finding.\n//  - Paragraph 2 (\"B\") becomes the adjacency-matrix finding.\n//  - Every existing list-item / traceback-line paragraph gets its run(s)\n//    (and the paragraph mark) bumped to 12pt (w:sz/w:szCs = 24 half-points).\n//  - A trailing bullet item (previously empty) gets new text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst newTextByIndex = {\n  0: \"Was able to simulate the PC.skel method for conditional independence.\",\n  1: \"Was able to create an output of an adjacency matrix using the MXM method of PC looking into how to apply the graphing of the causal learn package using the aforementioned adjacency matrix.\",\n  22: \"Current tech report pushed to iRCT repository.\",\n};\n\nfor (let i = 0; i < items.length; i++) {\n  const p = items[i];\n  const text = newTextByIndex[i];\n  if (text !== undefined) {\n    // Replaces (or, for the empty trailing bullet, inserts) the paragraph's\n    // run text while keeping the paragraph's own formatting/list membership.\n    p.insertText(text, \"Replace\");\n  }\n  // Bump the whole paragraph (paragraph mark + every run) to 12pt. Setting\n  // both `size` (w:sz) and `sizeBidirectional` (w:szCs) mirrors the rPr the\n  // diff adds to both <w:pPr> and every <w:r>.\n  p.font.size = 12;\n  p.font.sizeBidirectional = 12;\n}\n\nawait context.sync();\n", "ps1": "# Apply Garrett's Findings 1.25.23 meeting-notes update:\n#  - Paragraph 1 (\"A\") becomes the PC.skel finding.\n#  - Paragraph 2 (\"B\") becomes the adjacency-matrix finding.\n#  - Every existing list-item / traceback-line paragraph gets its run(s)\n#    (and the paragraph mark) bumped to 12pt (w:sz/w:szCs = 24 half-points).\n#  - The trailing bullet item (previously empty) gets new text.\n$doc = $word.ActiveDocument\n$paragraphs = $doc.Paragraphs\n\n$newText = @{\n    1  = \"Was able to simulate the PC.skel method for conditional independence.\"\n    2  = \"Was able to create an output of an adjacency matrix using the MXM method of PC looking into how to apply the graphing of the causal learn package using the aforementioned adjacency matrix.\"\n    23 = \"Current tech report pushed to iRCT repository.\"\n}\n\n$count = $paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paragraphs.Item($i)\n    $r = $p.Range\n\n    if ($newText.ContainsKey($i)) {\n        $r.Text = $newText[$i]\n    }\n\n    # Bump the whole paragraph (paragraph mark + every run) to 12pt.\n    $r.Font.Size = 12\n    $r.Font.SizeBi = 12\n}\n"}
